$wb = $excel.ActiveWorkbook

# --- Sheet 1: ProductSheet ---------------------------------------------
$ws1 = $wb.Worksheets.Item("ProductSheet")
$ws1.Range("B2").Value = "TEST-A-32"
$ws1.Range("B3").Value = "TEST-A-33"
$ws1.Range("B4").Value = "TEST-A-34"
$ws1.Range("B5").Value = "TEST-A-35"
$ws1.Range("B6").Value = "TEST-A-36"
$ws1.Range("B7").Value = "TEST-A-37"

# --- Sheet 2: editSheet --------------------------------------------------
$ws2 = $wb.Worksheets.Item("editSheet")
$ws2.Range("A2").Value = "TEST-A-32"
$ws2.Range("A3").Value = "TEST-A-33"
$ws2.Range("A4").Value = "TEST-A-34"
$ws2.Range("A5").Value = "TEST-A-35"
$ws2.Range("A6").Value = "TEST-A-36"
$ws2.Range("A7").Value = "TEST-A-37"

# --- Sheet 3: synchronization_sheet --------------------------------------
$ws3 = $wb.Worksheets.Item("synchronization_sheet")
$ws3.Range("B2").Value = "TEST-AS-31"
$ws3.Range("B3").Value = "TEST-AS-32"
$ws3.Range("B4").Value = "TEST-AS-33"
$ws3.Range("B5").Value = "TEST-AS-34"
$ws3.Range("B6").Value = "TEST-AS-35"
$ws3.Range("B7").Value = "TEST-AS-36"

# --- Selections / active sheet -------------------------------------------
[void]$ws2.Range("A2:A7").Select()
[void]$ws3.Range("B2:B7").Select()

[void]$ws1.Activate()
[void]$ws1.Range("B2:B7").Select()
